# "Sheet2" -> "구군" 이름 변경 및 활성 시트로 전환
$wb = $excel.ActiveWorkbook

$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Name = "구군"
$ws2.Activate()
